$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 3.2
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.6
$ws.Range("AC3").Value = 7.5
$ws.Range("AE3").Value = 17
$ws.Range("AM3").Value = 451
$ws.Range("AP3").Value = 23
$ws.Range("G4").Value = 2.3
$ws.Range("I4").Value = 3.7
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("U4").Value = 2.63
$ws.Range("V4").Value = 1.44
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 26
$ws.Range("AG4").Value = 7
$ws.Range("AX4").Value = 26
$ws.Range("AZ4").Value = 101
$ws.Range("O5").Value = 1.57
$ws.Range("P5").Value = 2.25
$ws.Range("G7").Value = 2.25
$ws.Range("H7").Value = 3
$ws.Range("K7").Value = 1.83
$ws.Range("L7").Value = 4.5
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6
$ws.Range("O7").Value = 1.57
$ws.Range("P7").Value = 2.25
$ws.Range("Q7").Value = 2.88
$ws.Range("R7").Value = 1.4
$ws.Range("U7").Value = 2.25
$ws.Range("V7").Value = 1.57
$ws.Range("AO7").Value = 15
$ws.Range("AX7").Value = 23
$ws.Range("G8").Value = 1.57
$ws.Range("H8").Value = 3.9
$ws.Range("I8").Value = 5.5
$ws.Range("J8").Value = 2.2
$ws.Range("K8").Value = 2.25
$ws.Range("L8").Value = 6
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.93
$ws.Range("U8").Value = 1.91
$ws.Range("V8").Value = 1.8
$ws.Range("X8").Value = 7
$ws.Range("Z8").Value = 11
$ws.Range("AD8").Value = 7.5
$ws.Range("AG8").Value = 15
$ws.Range("AH8").Value = 29
$ws.Range("AN8").Value = 3.5
$ws.Range("AO8").Value = 8
$ws.Range("AQ8").Value = 26
$ws.Range("AW8").Value = 7
$ws.Range("BB8").Value = 301
$ws.Range("G10").Value = 3.3
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 3
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 2.75
$ws.Range("Q10").Value = 2.3
$ws.Range("R10").Value = 1.6
$ws.Range("W10").Value = 8.5
$ws.Range("Y10").Value = 12
$ws.Range("AC10").Value = 7.5
$ws.Range("AH10").Value = 10
$ws.Range("AI10").Value = 9.5
$ws.Range("AU10").Value = 8.5
$ws.Range("G11").Value = 2.1
$ws.Range("I11").Value = 3.9
$ws.Range("J11").Value = 2.88
$ws.Range("L11").Value = 4.5
$ws.Range("W11").Value = 6
$ws.Range("X11").Value = 9
$ws.Range("AH11").Value = 19
$ws.Range("AI11").Value = 15
$ws.Range("AO11").Value = 12
$ws.Range("BA11").Value = 126
